# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund-level holdings) right before the
#   "总计" (totals) sheet.
# - Refresh the "总计" sheet so it gains a new top row summarizing 2022-Q1
#   (dates stay newest-first).
#
# The engine assigns a newly-created sheet's internal sheetId as
# (current max sheetId + 1). To land on the same ids the real edit used
# (2022-Q1 -> 4, 总计 -> 5, i.e. 总计 keeps growing past the id it used to
# hold) we delete the old 总计 sheet first (freeing id 4) and then add the
# two replacement sheets in order, so 2022-Q1 claims id 4 and the
# recreated 总计 claims id 5.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsTotalOld = $wb.Worksheets.Item("总计")
$wsTotalOld.Delete() | Out-Null

$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# Bold + thin-bordered + centered cell, matching the header/index-column
# look used throughout the workbook.
function Format-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Plain text cell with no special number formatting -- used for values
# that look numeric (e.g. "10.49") but must be stored as text, same as
# the source data. Forcing a text NumberFormat, assigning the value, then
# clearing formats keeps the stored type as text while leaving the cell
# with the workbook's default (unstyled) appearance.
function Set-PlainTextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# "2022-Q1" sheet: per-fund holdings for the quarter
# ---------------------------------------------------------------------

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $cell = $wsQ1.Cells.Item(1, 2 + $i)
    $cell.Value = $q1Headers[$i]
    Format-HeaderCell $cell
}

# code, name, scale, stock position, position pct, market value, rank
$q1Rows = @(
    @("004495", "博时量化平衡混合", "10.49", "38.32", "2.22", "0.2329", 2),
    @("012221", "瑞达行业轮动混合型证券投资基金A", "0.50", "84.76", "5.26", "0.0263", 5),
    @("516910", "富国中证现代物流交易型开放式指数证券投资基金", "0.42", "98.43", "4.82", "0.0202", 3),
    @("012222", "瑞达行业轮动混合型证券投资基金C", "0.17", "84.76", "5.26", "0.0089", 5)
)

for ($r = 0; $r -lt $q1Rows.Length; $r++) {
    $row = $q1Rows[$r]
    $excelRow = 2 + $r

    $idxCell = $wsQ1.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    Format-HeaderCell $idxCell

    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 2) $row[0]
    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 3) $row[1]
    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 4) $row[2]
    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 5) $row[3]
    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 6) $row[4]
    Set-PlainTextCell $wsQ1.Cells.Item($excelRow, 7) $row[5]
    $wsQ1.Cells.Item($excelRow, 8).Value = $row[6]
}

$wsQ1.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# "总计" sheet: one summary row per quarter, newest first
# ---------------------------------------------------------------------

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $cell = $wsTotal.Cells.Item(1, 2 + $i)
    $cell.Value = $totalHeaders[$i]
    Format-HeaderCell $cell
}

# date, count, market value
$totalRows = @(
    @("2022-Q1", 4, 0.29),
    @("2021-Q4", 2, 0.48),
    @("2021-Q3", 2, 0.03),
    @("2021-Q2", 1, 0.02)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = $totalRows[$r]
    $excelRow = 2 + $r

    $idxCell = $wsTotal.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    Format-HeaderCell $idxCell

    Set-PlainTextCell $wsTotal.Cells.Item($excelRow, 2) $row[0]
    $wsTotal.Cells.Item($excelRow, 3).Value = $row[1]
    $wsTotal.Cells.Item($excelRow, 4).Value = $row[2]
}

$wsTotal.Range("A1").Select() | Out-Null
